$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the "Statistiques globales" paragraph, then the two empty
# paragraphs that immediately follow it (a third empty paragraph follows
# and must stay untouched).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Statistiques globales*") {
        $anchor = $p
        break
    }
}

$firstEmpty = $anchor.Next()
$secondEmpty = $firstEmpty.Next()

# First empty paragraph -> "Stockage :" heading (same look as the other
# section headings, e.g. "Administration :").
$stockageXml = "<w:p $wNs>" +
    "<w:pPr><w:rPr><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr><w:t>Stockage :</w:t></w:r>" +
    "</w:p>"
[void]$firstEmpty.Range.InsertXML($stockageXml)

# Second empty paragraph -> three new bullet items reusing the existing
# "Paragraphedeliste" list (numId 3) already used above.
$item1 = "Base de données reliant les différents amis"
$item2 = "Base de données avec chaque utilisateur et les défis auxquels il participe"
$item3 = "Base de données avec les classements de chaque utilisateur (possibilité de regrouper par défi pour avoir les classements spécifiques"

$bulletXml = ""
foreach ($t in @($item1, $item2, $item3)) {
    $bulletXml += "<w:p $wNs>" +
        "<w:pPr><w:pStyle w:val='Paragraphedeliste'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='3'/></w:numPr></w:pPr>" +
        "<w:r><w:t>$t</w:t></w:r>" +
        "</w:p>"
}
[void]$secondEmpty.Range.InsertXML($bulletXml)
